# Apply the cryptos-list price/volume refresh described by the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.799.88'
$ws.Range("E2").Value = '  -0.58%  '
$ws.Range("D3").Value = '3.505.20'
$ws.Range("E3").Value = '  -1.53%  '
$ws.Range("E4").Value = '  -0.14%  '
# D5: plain-looking number ("606.81") must stay TEXT like the source cell,
# so force Text format before writing it, then restore the default style.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '606.81'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.30%  '
# D6: plain-looking number ("198.63") must stay TEXT like the source cell,
# so force Text format before writing it, then restore the default style.
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '198.63'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.76%  '
# D7: plain-looking number ("0.626") must stay TEXT like the source cell,
# so force Text format before writing it, then restore the default style.
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.626'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.99%  '
$ws.Range("E8").Value = '  -0.08%  '
# D9: plain-looking number ("0.211") must stay TEXT like the source cell,
# so force Text format before writing it, then restore the default style.
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.211'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.41%  '
# D10: plain-looking number ("0.657") must stay TEXT like the source cell,
# so force Text format before writing it, then restore the default style.
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.657'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.77%  '
# D11: plain-looking number ("54.20") must stay TEXT like the source cell,
# so force Text format before writing it, then restore the default style.
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '54.20'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.53%  '
# D12: plain-looking number ("0.0000305") must stay TEXT like the source cell,
# so force Text format before writing it, then restore the default style.
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000305'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.36%  '
# D13: plain-looking number ("9.63") must stay TEXT like the source cell,
# so force Text format before writing it, then restore the default style.
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.63'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.66%  '
$ws.Range("D14").Value = '4.068.84'
$ws.Range("E14").Value = '  -1.46%  '
# D15: plain-looking number ("597.54") must stay TEXT like the source cell,
# so force Text format before writing it, then restore the default style.
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '597.54'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.77%  '
$ws.Range("D16").Value = '69.978.13'
$ws.Range("E16").Value = '  -0.49%  '
# D17: plain-looking number ("19.00") must stay TEXT like the source cell,
# so force Text format before writing it, then restore the default style.
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '19.00'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.26%  '
# D18: plain-looking number ("12.67") must stay TEXT like the source cell,
# so force Text format before writing it, then restore the default style.
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.67'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.07%  '
$ws.Range("D19").Value = '3.518.15'
$ws.Range("E19").Value = '  -1.06%  '
$ws.Range("E20").Value = '  -0.16%  '
# D21: plain-looking number ("0.997") must stay TEXT like the source cell,
# so force Text format before writing it, then restore the default style.
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.997'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.48%  '
# D22: plain-looking number ("17.76") must stay TEXT like the source cell,
# so force Text format before writing it, then restore the default style.
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '17.76'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.19%  '
# D23: plain-looking number ("104.10") must stay TEXT like the source cell,
# so force Text format before writing it, then restore the default style.
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '104.10'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +10.79%  '
# D24: plain-looking number ("4.65") must stay TEXT like the source cell,
# so force Text format before writing it, then restore the default style.
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.65'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.27%  '
# D25: plain-looking number ("5.06") must stay TEXT like the source cell,
# so force Text format before writing it, then restore the default style.
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.06'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.89%  '
$ws.Range("E26").Value = '  +6.41%  '
# D27: plain-looking number ("10.99") must stay TEXT like the source cell,
# so force Text format before writing it, then restore the default style.
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.99'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.66%  '
# D28: plain-looking number ("9.82") must stay TEXT like the source cell,
# so force Text format before writing it, then restore the default style.
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.82'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.32%  '
# D29: plain-looking number ("33.85") must stay TEXT like the source cell,
# so force Text format before writing it, then restore the default style.
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.85'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.79%  '
# D30: plain-looking number ("4.58") must stay TEXT like the source cell,
# so force Text format before writing it, then restore the default style.
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.58'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +22.78%  '
# D31: plain-looking number ("7.22") must stay TEXT like the source cell,
# so force Text format before writing it, then restore the default style.
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.22'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.52%  '
$ws.Range("E32").Value = '  +4.48%  '
$ws.Range("E33").Value = '  +1.47%  '
# D34: plain-looking number ("63.84") must stay TEXT like the source cell,
# so force Text format before writing it, then restore the default style.
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '63.84'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.11%  '
$ws.Range("D35").Value = '3.721.40'
$ws.Range("E35").Value = '  +2.47%  '
# D36: plain-looking number ("522.04") must stay TEXT like the source cell,
# so force Text format before writing it, then restore the default style.
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '522.04'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.64%  '
$ws.Range("E37").Value = '  -0.24%  '
$ws.Range("D38").Value = '0.0₃0795'
$ws.Range("E38").Value = '  +1.67%  '
# D39: plain-looking number ("3.02") must stay TEXT like the source cell,
# so force Text format before writing it, then restore the default style.
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.02'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.98%  '
# D40: plain-looking number ("0.392") must stay TEXT like the source cell,
# so force Text format before writing it, then restore the default style.
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.392'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.87%  '
# D41: plain-looking number ("36.87") must stay TEXT like the source cell,
# so force Text format before writing it, then restore the default style.
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '36.87'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.47%  '
# D42: plain-looking number ("3.56") must stay TEXT like the source cell,
# so force Text format before writing it, then restore the default style.
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.56'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.41%  '
$ws.Range("E43").Value = '  -0.18%  '
# D44: plain-looking number ("0.0463") must stay TEXT like the source cell,
# so force Text format before writing it, then restore the default style.
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0463'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.48%  '
# D45: plain-looking number ("2.86") must stay TEXT like the source cell,
# so force Text format before writing it, then restore the default style.
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.86'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.56%  '
$ws.Range("E46").Value = '  +0.95%  '
$ws.Range("E47").Value = '  -4.60%  '
# D48: plain-looking number ("8.77") must stay TEXT like the source cell,
# so force Text format before writing it, then restore the default style.
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.77'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.48%  '
$ws.Range("E49").Value = '  +0.37%  '
# D50: plain-looking number ("132.22") must stay TEXT like the source cell,
# so force Text format before writing it, then restore the default style.
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '132.22'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.48%  '
# D51: plain-looking number ("0.000239") must stay TEXT like the source cell,
# so force Text format before writing it, then restore the default style.
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.000239'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.44%  '
